# BP-359 Bank excel statemenst upload
#
# Re-key the REF_NO column (B) as text (it was numeric before, and the new
# values collide per-row with the TXN_REF_NO prefix digit), tag the
# identifier columns (A/B/C) with a text format and the amount column (E)
# with a 2-decimal numeric format, move the active selection to C3, and set
# the sheet to print in portrait orientation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- REF_NO column (B) becomes text, with new per-row values ---------------
$ws.Range("B2:B6").NumberFormat = "@"
$ws.Range("B2").Value = "199999"
$ws.Range("B3").Value = "288888"
$ws.Range("B4").Value = "388888"
$ws.Range("B5").Value = "488888"
$ws.Range("B6").Value = "588888"

# --- A and C also tagged with the text format (values unchanged) -----------
$ws.Range("A2:A6").NumberFormat = "@"
$ws.Range("C2:C6").NumberFormat = "@"

# --- E (amount) gets a 2-decimal numeric format -----------------------------
$ws.Range("E2:E6").NumberFormat = "0.00"

# --- Selection moves to C3 ---------------------------------------------------
$ws.Range("C3").Select() | Out-Null

# --- Print orientation: portrait --------------------------------------------
$ws.PageSetup.Orientation = 1
